$d = $word.ActiveDocument

# Remove the existing "_GoBack" bookmark (it currently sits after "...УПЗ"
# near "...УПЗ-11"); it will be re-created at its new location below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Replace the title text "ВСТУП ДО СПЕЦІАЛЬНОСТІ" (which lives in two runs)
# with "МЕТОДОЛОГІЯ, МЕТОДИ І ЗАСОБИ УПРАВЛІННЯ ПРОЕКТАМИ".
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$found = $find.Execute("ВСТУП ДО СПЕЦІАЛЬНОСТІ", $false, $false, $false, $false, $false, $true, 1, $false, "МЕТОДОЛОГІЯ, МЕТОДИ І ЗАСОБИ УПРАВЛІННЯ ПРОЕКТАМИ", 2)

# Re-insert the "_GoBack" bookmark right after the freshly-typed text
# (i.e. where the cursor was left after the last edit), before the
# closing "»".
$find2 = $d.Content.Find
$find2.ClearFormatting()
$found2 = $find2.Execute("МЕТОДОЛОГІЯ, МЕТОДИ І ЗАСОБИ УПРАВЛІННЯ ПРОЕКТАМИ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endRng = $d.Content
$endRng.Start = $find2.Parent.End
$endRng.End = $find2.Parent.End
$d.Bookmarks.Add("_GoBack", $endRng)
